$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("login")
$wsRun   = $wb.Worksheets.Item("RUNMANAGER")

# --- "login" sheet: widen column A (testname) ---
$wsLogin.Columns.Item(1).ColumnWidth = 28.14

# --- "login" sheet: flip "execute" column from yes -> no on rows 2, 4 and 6 ---
$wsLogin.Range("B2").Value = "no"
$wsLogin.Range("B4").Value = "no"
$wsLogin.Range("B6").Value = "no"

# --- "login" sheet: touch up C5/D5 (username/password) so they match row formatting ---
$wsLogin.Range("C5").Value = "John Doe"
$wsLogin.Range("D5").Value = "ThisIsNotAPassword"

# --- "login" sheet: add the new "loginTest" row (row 7) ---
# Write E7 before A7 so the shared-string table gains "," (21) ahead of "loginTest" (22),
# matching how the strings were appended in the authored workbook.
$wsLogin.Range("E7").Value = ","
$wsLogin.Range("D7").Value = "ThisIsNotAPassword"
$wsLogin.Range("C7").Value = "John Doe"
$wsLogin.Range("B7").Value = "yes"
$wsLogin.Range("A7").Value = "loginTest"

# --- "login" sheet: selection / active cell now sits on B2 ---
$wsLogin.Range("B2").Select()

# --- "login" becomes the active (visible) tab; RUNMANAGER loses that state ---
$wsLogin.Activate()

$wb.Save()
